# Fixes datetime picker format, modal background issue, and destroys picker
# when necessary.
#
# The eonasdan (Bootstrap DateTime Picker) form sheet is updated so the
# picker's example rows carry separate date-only / time-only examples in
# addition to the existing combined date+time example, and the input
# attribute that used to target a single `data-field` is replaced with a
# dedicated `timeFormat` attribute (the old `min`/`max` input-attribute
# columns, which were never populated, are removed).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("survey")

# Drop the unused inputAttributes.min / inputAttributes.max columns.
$ws1.Range("O1:P1").EntireColumn.Delete()

# Row 2 (datetime / example1 / "Date and Time") is already in place.
# Add a date-only example and a time-only example as rows 3 and 4.
$ws1.Range("I3").Value = "example2"
$ws1.Range("I4").Value = "example3"
$ws1.Range("J3").Value = "Date"
$ws1.Range("J4").Value = "Time"
$ws1.Range("G3").Value = "date"
$ws1.Range("G4").Value = "time"

# The remaining inputAttributes.data-field column becomes
# inputAttributes.timeFormat, with format values for the new rows.
$ws1.Range("N1").Value = "inputAttributes.timeFormat"
$ws1.Range("N3").Value = "YYYY/DD/MM"
$ws1.Range("N4").Value = "HH:mm"

$ws1.PageSetup.Orientation = 1

# Make the survey sheet the active tab, with the newly added cell selected.
[void]$ws1.Activate()
[void]$ws1.Range("N4").Select()
